$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column C for rows 2-6 from
# 2023-10-09 (45208) to 2023-10-13 (45212), keeping the date formatting.
$ws.Range("C2:C6").Value = 45212
